$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, derived from a permutation of the existing rows'
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), O (Origen) and P (Precio $/Kg) data.
$updates = @{
    2  = @{ D = 44421; J = 20; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    3  = @{ D = 44446; J = 15; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 520 }
    7  = @{ D = 44425; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 560 }
    9  = @{ D = 44340; J = 25; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    10 = @{ D = 44418; J = 12; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 600 }
    11 = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 480 }
    12 = @{ D = 44432; J = 15; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 }
    13 = @{ D = 44449; J = 30; K = 16000; L = 16000; M = 16000; O = "Provincia de Limarí"; P = 640 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 10).Value = $vals.J
    $ws.Cells.Item($row, 11).Value = $vals.K
    $ws.Cells.Item($row, 12).Value = $vals.L
    $ws.Cells.Item($row, 13).Value = $vals.M
    $ws.Cells.Item($row, 15).Value = $vals.O
    $ws.Cells.Item($row, 16).Value = $vals.P
}
